{"js": "// Replace each three-digit-by-one-digit multiplication expression\n// with its updated counterpart, per the commit diff.\nconst replacements = [\n  [\"845\u00d77=5915\", \"162\u00d76=972\"],\n  [\"145\u00d72=290\", \"756\u00d74=3024\"],\n  [\"984\u00d74=3936\", \"759\u00d78=6072\"],\n  [\"839\u00d73=2517\", \"900\u00d74=3600\"],\n  [\"551\u00d79=4959\", \"707\u00d79=6363\"],\n  [\"790\u00d73=2370\", \"239\u00d77=1673\"],\n  [\"523\u00d75=2615\", \"485\u00d75=2425\"],\n  [\"645\u00d77=4515\", \"939\u00d74=3756\"],\n  [\"337\u00d79=3033\", \"209\u00d79=1881\"],\n  [\"888\u00d78=7104\", \"242\u00d76=1452\"],\n  [\"719\u00d73=2157\", \"584\u00d73=1752\"],\n  [\"593\u00d76=3558\", \"890\u00d77=6230\"],\n  [\"828\u00d72=1656\", \"594\u00d79=5346\"],\n  [\"810\u00d79=7290\", \"341\u00d75=1705\"],\n  [\"891\u00d78=7128\", \"485\u00d77=3395\"],\n  [\"930\u00d76=5580\", \"492\u00d78=3936\"],\n  [\"822\u00d72=1644\", \"898\u00d79=8082\"],\n  [\"656\u00d74=2624\", \"134\u00d74=536\"],\n  [\"885\u00d79=7965\", \"670\u00d78=5360\"],\n  [\"159\u00d79=1431\", \"676\u00d76=4056\"],\n  [\"749\u00d75=3745\", \"114\u00d76=684\"],\n  [\"198\u00d75=990\", \"945\u00d75=4725\"],\n  [\"506\u00d79=4554\", \"352\u00d74=1408\"],\n  [\"982\u00d72=1964\", \"685\u00d73=2055\"],\n  [\"897\u00d79=8073\", \"998\u00d74=3992\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, 'Replace');\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update each three-digit-by-one-digit multiplication expression\n# with its updated counterpart, per the commit diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@(\"845\u00d77=5915\", \"162\u00d76=972\")\n    ,@(\"145\u00d72=290\", \"756\u00d74=3024\")\n    ,@(\"984\u00d74=3936\", \"759\u00d78=6072\")\n    ,@(\"839\u00d73=2517\", \"900\u00d74=3600\")\n    ,@(\"551\u00d79=4959\", \"707\u00d79=6363\")\n    ,@(\"790\u00d73=2370\", \"239\u00d77=1673\")\n    ,@(\"523\u00d75=2615\", \"485\u00d75=2425\")\n    ,@(\"645\u00d77=4515\", \"939\u00d74=3756\")\n    ,@(\"337\u00d79=3033\", \"209\u00d79=1881\")\n    ,@(\"888\u00d78=7104\", \"242\u00d76=1452\")\n    ,@(\"719\u00d73=2157\", \"584\u00d73=1752\")\n    ,@(\"593\u00d76=3558\", \"890\u00d77=6230\")\n    ,@(\"828\u00d72=1656\", \"594\u00d79=5346\")\n    ,@(\"810\u00d79=7290\", \"341\u00d75=1705\")\n    ,@(\"891\u00d78=7128\", \"485\u00d77=3395\")\n    ,@(\"930\u00d76=5580\", \"492\u00d78=3936\")\n    ,@(\"822\u00d72=1644\", \"898\u00d79=8082\")\n    ,@(\"656\u00d74=2624\", \"134\u00d74=536\")\n    ,@(\"885\u00d79=7965\", \"670\u00d78=5360\")\n    ,@(\"159\u00d79=1431\", \"676\u00d76=4056\")\n    ,@(\"749\u00d75=3745\", \"114\u00d76=684\")\n    ,@(\"198\u00d75=990\", \"945\u00d75=4725\")\n    ,@(\"506\u00d79=4554\", \"352\u00d74=1408\")\n    ,@(\"982\u00d72=1964\", \"685\u00d73=2055\")\n    ,@(\"897\u00d79=8073\", \"998\u00d74=3992\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
